$d = $word.ActiveDocument

$d.Content.Find.Execute(
    "Ημερομηνίες παρατήρησης για τον αστερισμό του Αστερισμός του Ηρακλή 2022: 13-22 Ιουνίου, 12-21 Ιουλίου, 10-19 Αυγούστου",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "2022 Ημερομηνίες παρατήρησης για τον αστερισμό του Αστερισμός του Ηρακλή: 13-22 Ιουνίου, 12-21 Ιουλίου, 10-19 Αυγούστου",
    2
)
